$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 61666.668
$ws.Range("J3").Value = 61666.668
$ws.Range("L3").Value = 61666.668
$ws.Range("N3").Value = -61894.668
$ws.Range("H33").Value = 15152024
$ws.Range("I33").Value = 22222794
$ws.Range("J33").Value = 373.14285
$ws.Range("K33").Value = 22222794
$ws.Range("L33").Value = 373.14285
$ws.Range("M33").Value = -22222565
$ws.Range("N33").Value = -831.14285
$ws.Range("H100").Value = 2085
$ws.Range("I100").Value = 1398.75
$ws.Range("K100").Value = 1398.75
$ws.Range("M100").Value = -857.75
$ws.Range("H102").Value = 61666.668
$ws.Range("J102").Value = 61666.668
$ws.Range("L102").Value = 61666.668
$ws.Range("N102").Value = -68156.668
$ws.Range("H138").Value = 3080.024
$ws.Range("I138").Value = 2543.76
$ws.Range("J138").Value = 3307.2542
$ws.Range("K138").Value = 7631.280000000001
$ws.Range("L138").Value = 9921.7626
$ws.Range("M138").Value = -2491.280000000001
$ws.Range("N138").Value = -20201.7626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25357.988
$ws.Range("I32").Value = 25454.49
$ws.Range("J32").Value = 25153.846
$ws.Range("K32").Value = 25454.49
$ws.Range("L32").Value = 25153.846
$ws.Range("M32").Value = -25167.49
$ws.Range("N32").Value = -25727.846
$ws.Range("H61").Value = 1433.125
$ws.Range("I61").Value = 1176.9783
$ws.Range("K61").Value = 1176.9783
$ws.Range("M61").Value = -964.9783
$ws.Range("H110").Value = 2504.75
$ws.Range("I110").Value = 2640.1177
$ws.Range("J110").Value = 1737.6666
$ws.Range("K110").Value = 2640.1177
$ws.Range("L110").Value = 1737.6666
$ws.Range("M110").Value = -595.1176999999998
$ws.Range("N110").Value = -5827.6666
$ws.Range("H136").Value = 1433.125
$ws.Range("I136").Value = 1176.9783
$ws.Range("K136").Value = 3530.9349
$ws.Range("M136").Value = -980.9349000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 356.26923
$ws.Range("I94").Value = 367.25
$ws.Range("J94").Value = 224.5
$ws.Range("K94").Value = 367.25
$ws.Range("L94").Value = 224.5
$ws.Range("M94").Value = 83.75
$ws.Range("N94").Value = -1126.5
$ws.Range("H134").Value = 4470.8076
$ws.Range("I134").Value = 6295.423
$ws.Range("J134").Value = 2646.1924
$ws.Range("K134").Value = 18886.269
$ws.Range("L134").Value = 7938.5772
$ws.Range("M134").Value = -16351.269
$ws.Range("N134").Value = -13008.5772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I68").Value = 499
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1497
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -686
$ws.Range("N68").ClearContents()
$ws.Range("I71").Value = 499
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 4491
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -435
$ws.Range("N71").ClearContents()
$ws.Range("H131").Value = 2313.8484
$ws.Range("I131").Value = 50580
$ws.Range("J131").Value = 805.53125
$ws.Range("K131").Value = 151740
$ws.Range("L131").Value = 2416.59375
$ws.Range("M131").Value = -146700
$ws.Range("N131").Value = -12496.59375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 398.81818
$ws.Range("I107").Value = 459.33334
$ws.Range("J107").Value = 126.5
$ws.Range("K107").Value = 459.33334
$ws.Range("L107").Value = 126.5
$ws.Range("M107").Value = 1460.66666
$ws.Range("N107").Value = -3966.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15874336
$ws.Range("I61").Value = 1181
$ws.Range("J61").Value = 33334806
$ws.Range("K61").Value = 1181
$ws.Range("L61").Value = 33334806
$ws.Range("M61").Value = -979
$ws.Range("N61").Value = -33335210
$ws.Range("H100").Value = 111279450
$ws.Range("I100").Value = 111279450
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 111279450
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -111278909
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 15874336
$ws.Range("I113").Value = 1181
$ws.Range("J113").Value = 33334806
$ws.Range("K113").Value = 1181
$ws.Range("L113").Value = 33334806
$ws.Range("M113").Value = 989
$ws.Range("N113").Value = -33339146
$ws.Range("H122").Value = 4818.5
$ws.Range("I122").Value = 6040.1
$ws.Range("J122").Value = 2782.5
$ws.Range("K122").Value = 18120.3
$ws.Range("L122").Value = 8347.5
$ws.Range("M122").Value = -15670.3
$ws.Range("N122").Value = -13247.5
$ws.Range("H136").Value = 4639.5713
$ws.Range("I136").Value = 5225.2256
$ws.Range("J136").Value = 2989.0908
$ws.Range("K136").Value = 15675.6768
$ws.Range("L136").Value = 8967.2724
$ws.Range("M136").Value = -13125.6768
$ws.Range("N136").Value = -14067.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13131.091
$ws.Range("J41").Value = 12011.111
$ws.Range("L41").Value = 12011.111
$ws.Range("N41").Value = -12791.111
$ws.Range("H45").Value = 9288.625
$ws.Range("J45").Value = 9288.625
$ws.Range("L45").Value = 9288.625
$ws.Range("N45").Value = -10270.625
$ws.Range("H70").Value = 21966.666
$ws.Range("I70").Value = 7000
$ws.Range("J70").Value = 29450
$ws.Range("K70").Value = 7000
$ws.Range("L70").Value = 29450
$ws.Range("M70").Value = -6685
$ws.Range("N70").Value = -30080
$ws.Range("H73").Value = 21966.666
$ws.Range("I73").Value = 7000
$ws.Range("J73").Value = 29450
$ws.Range("K73").Value = 7000
$ws.Range("L73").Value = 29450
$ws.Range("M73").Value = -5908
$ws.Range("N73").Value = -31634
$ws.Range("H86").Value = 20980.5
$ws.Range("J86").Value = 20980.5
$ws.Range("L86").Value = 20980.5
$ws.Range("N86").Value = -23226.5
$ws.Range("H87").Value = 76000
$ws.Range("J87").Value = 76000
$ws.Range("L87").Value = 76000
$ws.Range("N87").Value = -78496
$ws.Range("H89").Value = 20980.5
$ws.Range("J89").Value = 20980.5
$ws.Range("L89").Value = 104902.5
$ws.Range("N89").Value = -116134.5
$ws.Range("H90").Value = 76000
$ws.Range("J90").Value = 76000
$ws.Range("L90").Value = 228000
$ws.Range("N90").Value = -240480
$ws.Range("H126").Value = 41674096
$ws.Range("I126").Value = 55563956
$ws.Range("J126").Value = 4510.6665
$ws.Range("K126").Value = 166691868
$ws.Range("L126").Value = 13531.9995
$ws.Range("M126").Value = -166689398
$ws.Range("N126").Value = -18471.9995
